# fix PlainExcelGenerator for DateOnly and TimeOnly
# Adds a new "Hour" sample column (I) to the plain-excel template, showing
# a TimeOnly-style value formatted as a time-of-day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the "Hour" column, next to the existing "Percentage" header.
$ws.Range("I2").Value = "Hour"

# Sample value for the new column, formatted like a TimeOnly value (12:00 PM).
$ws.Range("I3").Value = 0.5
$ws.Range("I3").NumberFormat = "[$-F400]h:mm:ss am/pm"

# Keep selection where the new content is, matching the updated template.
$ws.Range("I3").Select() | Out-Null
